$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("electricite", 45, 45536, "Cosy Appart - Plage à 2min - Casino à 1min"),
    @("copro", 70, 45536, "Cosy Appart - Plage à 2min - Casino à 1min"),
    @("box ", 30, 45536, "Cosy Appart - Plage à 2min - Casino à 1min"),
    @("credit", 880, 45536, "Cosy Appart - Plage à 2min - Casino à 1min"),
    @("samantha", 120, 45536, "Cosy Appart - Plage à 2min - Casino à 1min")
)

$row = 12
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 3).NumberFormat = "mmm-yy"
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$ws.Range("I24").Select() | Out-Null
